$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "22.477.26"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3
$ws.Range("D3").Value = "1.574.06"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.43"
$ws.Range("E6").Value = "  +0.13%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3750"
$ws.Range("E7").Value = "  -0.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.99"
$ws.Range("E8").Value = "  +0.17%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3407"
$ws.Range("E9").Value = "  -0.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.150"
$ws.Range("E10").Value = "  -1.12%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07587"
$ws.Range("E11").Value = "  -1.01%  "

# Row 12
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.40"
$ws.Range("E13").Value = "  +0.52%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.001"
$ws.Range("E14").Value = "  +0.41%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.972"

# Row 16
$ws.Range("D16").Value = "1.573.98"
$ws.Range("E16").Value = "  -0.14%  "

# Row 17
$ws.Range("E17").Value = "  -0.95%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.25"
$ws.Range("E18").Value = "  +0.72%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06744"
$ws.Range("E19").Value = "  -0.16%  "

# Row 20
$ws.Range("E20").Value = "  -0.17%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.283"
$ws.Range("E21").Value = "  +0.80%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.45"
$ws.Range("E22").Value = "  -1.87%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.19"
$ws.Range("E23").Value = "  +1.42%  "

# Row 24
$ws.Range("D24").Value = "22.468.15"
$ws.Range("E24").Value = "  +0.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.607"
$ws.Range("E26").Value = "  -5.54%  "

# Row 27
$ws.Range("E27").Value = "  -0.55%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.82"
$ws.Range("E28").Value = "  +2.41%  "

# Row 29
$ws.Range("E29").Value = "  -1.20%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.21"
$ws.Range("E30").Value = "  -0.02%  "

# Row 31
$ws.Range("D31").Value = "1.749.12"
$ws.Range("E31").Value = "  -0.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.042"
$ws.Range("E32").Value = "  +2.71%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.159"
$ws.Range("E33").Value = "  -0.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.978"
$ws.Range("E34").Value = "  -2.14%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.891"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08449"
$ws.Range("E36").Value = "  -0.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.388"
$ws.Range("E37").Value = "  +4.40%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02473"
$ws.Range("E38").Value = "  -3.45%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2298"
$ws.Range("E39").Value = "  -0.92%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06558"
$ws.Range("E40").Value = "  +0.30%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.515"
$ws.Range("E41").Value = "  +0.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.42"
$ws.Range("E42").Value = "  -1.74%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6302"
$ws.Range("E43").Value = "  -2.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.17"
$ws.Range("E44").Value = "  -0.07%  "

# Row 45
$ws.Range("E45").Value = "  -0.17%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.819"
$ws.Range("E46").Value = "  +0.69%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5889"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.100"

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "130.33"
$ws.Range("E49").Value = "  +3.57%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.230"
$ws.Range("E50").Value = "  -5.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07345"
$ws.Range("E51").Value = "  +0.01%  "
